# remove html br line breaks from case rule source
#
# The "rules.csv" sheet stores each rule's description/comment text as a
# shared string. A number of those strings contain literal "<br/>" markers
# (leftover HTML line breaks) mixed in with real newlines. Strip the stray
# "<br/>" markers wherever they show up in the used range, leaving
# everything else (values, numbers, formatting) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
  for ($c = 1; $c -le $colCount; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $text = [string]$cell.Text

    if ($text -ne $null -and $text.Contains("<br/>")) {
      $cleaned = $text.Replace("<br/>", "")
      $cell.Value = $cleaned
    }
  }
}
